$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06440233333333333
$ws.Range("H2").Value = 0.193207
$ws.Range("I2").Value = 0.03647206354366116
$ws.Range("J2").Value = 0.03647206354366116
$ws.Range("M2").Value = 29.52617166666667
$ws.Range("N2").Value = 88.57851500000001
$ws.Range("O2").Value = 0.3218391660320701
$ws.Range("P2").Value = 0.3218391660320701
$ws.Range("Q2").Value = 1.901554349733889
$ws.Range("R2").Value = 17.113989147605
$ws.Range("S2").Value = 0.01173813851436057
$ws.Range("T2").Value = 0.01173813851436057

$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06440233333333333
$ws.Range("H3").Value = 0.193207
$ws.Range("I3").Value = 0.03647206354366116
$ws.Range("J3").Value = 0.03647206354366116
$ws.Range("O3").Value = 0.4328989896002822
$ws.Range("P3").Value = 0.4328989896002822
$ws.Range("Q3").Value = 2.557740149587
$ws.Range("R3").Value = 23.019661346283
$ws.Range("S3").Value = 0.0157887194566882
$ws.Range("T3").Value = 0.0157887194566882

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06440233333333333
$ws.Range("H4").Value = 0.193207
$ws.Range("I4").Value = 0.03647206354366116
$ws.Range("J4").Value = 0.03647206354366116
$ws.Range("M4").Value = 22.50081433333333
$ws.Range("N4").Value = 67.502443
$ws.Range("O4").Value = 0.2452618443676477
$ws.Range("P4").Value = 0.2452618443676476
$ws.Range("Q4").Value = 1.449104944966778
$ws.Range("R4").Value = 13.041944504701
$ws.Range("S4").Value = 0.008945205572612379
$ws.Range("T4").Value = 0.008945205572612377

$ws.Range("I5").Value = 0.8194013021867156
$ws.Range("J5").Value = 0.8194013021867155
$ws.Range("M5").Value = 29.52617166666667
$ws.Range("N5").Value = 88.57851500000001
$ws.Range("O5").Value = 0.3218391660320701
$ws.Range("P5").Value = 0.3218391660320701
$ws.Range("Q5").Value = 42.72135873215667
$ws.Range("R5").Value = 384.49222858941
$ws.Range("S5").Value = 0.2637154317413648
$ws.Range("T5").Value = 0.2637154317413647

$ws.Range("I6").Value = 0.8194013021867156
$ws.Range("J6").Value = 0.8194013021867155
$ws.Range("O6").Value = 0.4328989896002822
$ws.Range("P6").Value = 0.4328989896002822
$ws.Range("S6").Value = 0.3547179957937847
$ws.Range("T6").Value = 0.3547179957937847

$ws.Range("I7").Value = 0.8194013021867156
$ws.Range("J7").Value = 0.8194013021867155
$ws.Range("M7").Value = 22.50081433333333
$ws.Range("N7").Value = 67.502443
$ws.Range("O7").Value = 0.2452618443676477
$ws.Range("P7").Value = 0.2452618443676476
$ws.Range("Q7").Value = 32.55638325727134
$ws.Range("R7").Value = 293.007449315442
$ws.Range("S7").Value = 0.2009678746515661
$ws.Range("T7").Value = 0.200967874651566

$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.2544986666666667
$ws.Range("H8").Value = 0.763496
$ws.Range("I8").Value = 0.1441266342696234
$ws.Range("J8").Value = 0.1441266342696234
$ws.Range("M8").Value = 29.52617166666667
$ws.Range("N8").Value = 88.57851500000001
$ws.Range("O8").Value = 0.3218391660320701
$ws.Range("P8").Value = 0.3218391660320701
$ws.Range("Q8").Value = 7.514371320937778
$ws.Range("R8").Value = 67.62934188844001
$ws.Range("S8").Value = 0.04638559577634475
$ws.Range("T8").Value = 0.04638559577634475

$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.2544986666666667
$ws.Range("H9").Value = 0.763496
$ws.Range("I9").Value = 0.1441266342696234
$ws.Range("J9").Value = 0.1441266342696234
$ws.Range("O9").Value = 0.4328989896002822
$ws.Range("P9").Value = 0.4328989896002822
$ws.Range("Q9").Value = 10.107420400136
$ws.Range("R9").Value = 90.96678360122399
$ws.Range("S9").Value = 0.06239227434980936
$ws.Range("T9").Value = 0.06239227434980936

$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.2544986666666667
$ws.Range("H10").Value = 0.763496
$ws.Range("I10").Value = 0.1441266342696234
$ws.Range("J10").Value = 0.1441266342696234
$ws.Range("M10").Value = 22.50081433333333
$ws.Range("N10").Value = 67.502443
$ws.Range("O10").Value = 0.2452618443676477
$ws.Range("P10").Value = 0.2452618443676476
$ws.Range("Q10").Value = 5.726427246747556
$ws.Range("R10").Value = 51.53784522072799
$ws.Range("S10").Value = 0.03534876414346924
$ws.Range("T10").Value = 0.03534876414346923

